$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates (August 2025 rows) ---
# Row 3: Dia 4 -> total_venda 30338.81 -> 30345.31
$ws.Cells.Item(3,2).Value = 30345.31

# Row 7: Dia 8 -> total_venda 13507.96 -> 13698.11
$ws.Cells.Item(7,2).Value = 13698.11

# --- Insert a new row (Dia 11, Ago/2025) right after the existing August rows.
#     This single insert shifts every following row down by one, which is all
#     that's needed to line everything up with the rest of the data through
#     to the end of the sheet (no further row insertion is required). ---
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8,1).Value = 11
$ws.Cells.Item(8,2).Value = 11901.05
$ws.Cells.Item(8,3).Value = 8
$ws.Cells.Item(8,4).Value = 2025
$ws.Cells.Item(8,5).Value = "08/2025"

# After the insertion, the row that used to be "Dia 14 / Jul 2025" (old row 17)
# is now row 18; its total_venda also changed: 18544.3 -> 22074.6
$ws.Cells.Item(18,2).Value = 22074.6
